$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy cell formatting (number formats/styles) from the old D column (now E) into the new D column
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Match the new column D's width to the other data columns (E:K)
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(8)

# Populate the newly inserted column D with the new fiscal-year data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 96500
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -400
$ws.Range("D17").Value = 15100
$ws.Range("D18").Value = 81400
$ws.Range("D20").Value = -36100
$ws.Range("D21").Value = 50100
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 45300
$ws.Range("D24").Value = 9000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 36300
$ws.Range("D27").Value = 36300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 36100
$ws.Range("D33").Value = 36300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 36300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 56500
$ws.Range("D42").Value = 43200
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 30400
$ws.Range("D49").Value = 23700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2988300
$ws.Range("D57").Value = "NA"
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 20000
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2718800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 29300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 269600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 36300
$ws.Range("D83").Value = 4800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 42000
$ws.Range("D91").Value = -5100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -227800
$ws.Range("D96").Value = -14400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 197200
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 11400
